# repull data, push all data, mean calculation
# Update the dSF (column F) values on Sheet1 to match the re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "F4"  = 0
    "F8"  = 0
    "F9"  = 1
    "F10" = 5
    "F13" = 0
    "F16" = -2
    "F19" = -1
    "F26" = -1
    "F27" = 0
    "F30" = 1
    "F34" = 0
    "F35" = -3
    "F36" = -9
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
